$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 1937.8823
$ws.Cells.Item(8, 9).Value = 74
$ws.Cells.Item(8, 11).Value = 222
$ws.Cells.Item(8, 13).Value = -83
$ws.Cells.Item(70, 8).Value = 3699.2
$ws.Cells.Item(70, 9).Value = 1500
$ws.Cells.Item(70, 10).Value = 4249
$ws.Cells.Item(70, 11).Value = 4500
$ws.Cells.Item(70, 12).Value = 12747
$ws.Cells.Item(70, 13).Value = -4230
$ws.Cells.Item(70, 14).Value = -13287
$ws.Cells.Item(73, 8).Value = 3699.2
$ws.Cells.Item(73, 9).Value = 1500
$ws.Cells.Item(73, 10).Value = 4249
$ws.Cells.Item(73, 11).Value = 4500
$ws.Cells.Item(73, 12).Value = 12747
$ws.Cells.Item(73, 13).Value = -3564
$ws.Cells.Item(73, 14).Value = -14619
$ws.Cells.Item(98, 8).Value = 902.86957
$ws.Cells.Item(98, 10).Value = 1166.3334
$ws.Cells.Item(98, 12).Value = 1166.3334
$ws.Cells.Item(98, 14).Value = -4162.3334
$ws.Cells.Item(106, 8).Value = 17800
$ws.Cells.Item(106, 9).Value = 14500
$ws.Cells.Item(106, 11).Value = 14500
$ws.Cells.Item(106, 13).Value = -13869
$ws.Cells.Item(112, 8).Value = 1355.7
$ws.Cells.Item(112, 9).Value = 1489.8334
$ws.Cells.Item(112, 10).Value = 1298.2142
$ws.Cells.Item(112, 11).Value = 4469.5002
$ws.Cells.Item(112, 12).Value = 3894.6426
$ws.Cells.Item(112, 13).Value = -3361.5002
$ws.Cells.Item(112, 14).Value = -6110.642599999999
$ws.Cells.Item(122, 8).Value = 902.86957
$ws.Cells.Item(122, 10).Value = 1166.3334
$ws.Cells.Item(122, 12).Value = 3499.0002
$ws.Cells.Item(122, 14).Value = -8399.0002
$ws.Cells.Item(125, 8).Value = 4519.5
$ws.Cells.Item(125, 9).Value = 4040.75
$ws.Cells.Item(125, 11).Value = 36366.75
$ws.Cells.Item(125, 13).Value = -33906.75
$ws.Cells.Item(137, 8).Value = 6457348.5
$ws.Cells.Item(137, 10).Value = 10143.467
$ws.Cells.Item(137, 12).Value = 30430.401
$ws.Cells.Item(137, 14).Value = -35530.401
$ws.Cells.Item(138, 8).Value = 7276.757
$ws.Cells.Item(138, 9).Value = 4973.2
$ws.Cells.Item(138, 10).Value = 8129.926
$ws.Cells.Item(138, 11).Value = 14919.6
$ws.Cells.Item(138, 12).Value = 24389.778
$ws.Cells.Item(138, 13).Value = -9779.599999999999
$ws.Cells.Item(138, 14).Value = -34669.77800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 489314.53
$ws.Cells.Item(32, 9).Value = 535123.7
$ws.Cells.Item(32, 11).Value = 535123.7
$ws.Cells.Item(32, 13).Value = -534836.7
$ws.Cells.Item(61, 8).Value = 14840865
$ws.Cells.Item(61, 9).Value = 5686027
$ws.Cells.Item(61, 11).Value = 5686027
$ws.Cells.Item(61, 13).Value = -5685815
$ws.Cells.Item(80, 8).Value = 41000
$ws.Cells.Item(80, 10).Value = 42000
$ws.Cells.Item(80, 12).Value = 42000
$ws.Cells.Item(80, 14).Value = -43996
$ws.Cells.Item(83, 8).Value = 41000
$ws.Cells.Item(83, 10).Value = 42000
$ws.Cells.Item(83, 12).Value = 126000
$ws.Cells.Item(83, 14).Value = -135984
$ws.Cells.Item(134, 8).Value = 74666.336
$ws.Cells.Item(134, 10).Value = 74666.336
$ws.Cells.Item(134, 12).Value = 74666.336
$ws.Cells.Item(134, 14).Value = -84806.336
$ws.Cells.Item(136, 8).Value = 14840865
$ws.Cells.Item(136, 9).Value = 5686027
$ws.Cells.Item(136, 11).Value = 17058081
$ws.Cells.Item(136, 13).Value = -17055531

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 416041.7
$ws.Cells.Item(20, 9).Value = 758052.4399999999
$ws.Cells.Item(20, 11).Value = 758052.4399999999
$ws.Cells.Item(20, 13).Value = -757805.4399999999
$ws.Cells.Item(82, 8).Value = 22270
$ws.Cells.Item(82, 9).Value = 10738.571
$ws.Cells.Item(82, 10).Value = 49176.668
$ws.Cells.Item(82, 11).Value = 10738.571
$ws.Cells.Item(82, 12).Value = 49176.668
$ws.Cells.Item(82, 13).Value = -10355.571
$ws.Cells.Item(82, 14).Value = -49942.668
$ws.Cells.Item(85, 8).Value = 22270
$ws.Cells.Item(85, 9).Value = 10738.571
$ws.Cells.Item(85, 10).Value = 49176.668
$ws.Cells.Item(85, 11).Value = 10738.571
$ws.Cells.Item(85, 12).Value = 49176.668
$ws.Cells.Item(85, 13).Value = -9412.571
$ws.Cells.Item(85, 14).Value = -51828.668
$ws.Cells.Item(86, 8).Value = 5504
$ws.Cells.Item(86, 9).Value = 5471.4517
$ws.Cells.Item(86, 10).Value = 5588.0835
$ws.Cells.Item(86, 11).Value = 5471.4517
$ws.Cells.Item(86, 12).Value = 5588.0835
$ws.Cells.Item(86, 13).Value = -4348.4517
$ws.Cells.Item(86, 14).Value = -7834.0835
$ws.Cells.Item(89, 8).Value = 5504
$ws.Cells.Item(89, 9).Value = 5471.4517
$ws.Cells.Item(89, 10).Value = 5588.0835
$ws.Cells.Item(89, 11).Value = 27357.2585
$ws.Cells.Item(89, 12).Value = 27940.4175
$ws.Cells.Item(89, 13).Value = -21741.2585
$ws.Cells.Item(89, 14).Value = -39172.4175
$ws.Cells.Item(109, 8).Value = 70000
$ws.Cells.Item(109, 10).Value = 70000
$ws.Cells.Item(109, 12).Value = 70000
$ws.Cells.Item(109, 14).Value = -72774

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 925.2222
$ws.Cells.Item(22, 10).Value = 1749.6666
$ws.Cells.Item(22, 12).Value = 1749.6666
$ws.Cells.Item(22, 14).Value = -2449.6666
$ws.Cells.Item(31, 8).Value = 500719.7
$ws.Cells.Item(31, 10).Value = 5057.2646
$ws.Cells.Item(31, 12).Value = 5057.2646
$ws.Cells.Item(31, 14).Value = -5647.2646
$ws.Cells.Item(34, 8).Value = 500719.7
$ws.Cells.Item(34, 10).Value = 5057.2646
$ws.Cells.Item(34, 12).Value = 5057.2646
$ws.Cells.Item(34, 14).Value = -5461.2646
$ws.Cells.Item(58, 8).Value = 13138231
$ws.Cells.Item(58, 9).Value = 16670592
$ws.Cells.Item(58, 10).Value = 7250962
$ws.Cells.Item(58, 11).Value = 16670592
$ws.Cells.Item(58, 12).Value = 7250962
$ws.Cells.Item(58, 13).Value = -16670389
$ws.Cells.Item(58, 14).Value = -7251368
$ws.Cells.Item(105, 8).Value = 12508.958
$ws.Cells.Item(105, 9).Value = 11153.096
$ws.Cells.Item(105, 10).Value = 22000
$ws.Cells.Item(105, 11).Value = 11153.096
$ws.Cells.Item(105, 12).Value = 22000
$ws.Cells.Item(105, 13).Value = -9406.096
$ws.Cells.Item(105, 14).Value = -25494
$ws.Cells.Item(136, 8).Value = 13138231
$ws.Cells.Item(136, 9).Value = 16670592
$ws.Cells.Item(136, 10).Value = 7250962
$ws.Cells.Item(136, 11).Value = 50011776
$ws.Cells.Item(136, 12).Value = 21752886
$ws.Cells.Item(136, 13).Value = -50009226
$ws.Cells.Item(136, 14).Value = -21757986

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 71508.64
$ws.Cells.Item(2, 10).Value = 108.333336
$ws.Cells.Item(2, 12).Value = 108.333336
$ws.Cells.Item(2, 14).Value = -334.333336
$ws.Cells.Item(80, 8).Value = 5000
$ws.Cells.Item(80, 9).Value = 2500
$ws.Cells.Item(80, 10).Value = 10000
$ws.Cells.Item(80, 11).Value = 2500
$ws.Cells.Item(80, 12).Value = 10000
$ws.Cells.Item(80, 13).Value = -1502
$ws.Cells.Item(80, 14).Value = -11996
$ws.Cells.Item(83, 8).Value = 5000
$ws.Cells.Item(83, 9).Value = 2500
$ws.Cells.Item(83, 10).Value = 10000
$ws.Cells.Item(83, 11).Value = 12500
$ws.Cells.Item(83, 12).Value = 50000
$ws.Cells.Item(83, 13).Value = -7508
$ws.Cells.Item(83, 14).Value = -59984
$ws.Cells.Item(102, 8).Value = 1715.7241
$ws.Cells.Item(102, 9).Value = 1706.2174
$ws.Cells.Item(102, 11).Value = 1706.2174
$ws.Cells.Item(102, 13).Value = -84.2174

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 6002
$ws.Cells.Item(68, 9).Value = 6002
$ws.Cells.Item(68, 11).Value = 6002
$ws.Cells.Item(68, 13).Value = -5253
$ws.Cells.Item(71, 8).Value = 6002
$ws.Cells.Item(71, 9).Value = 6002
$ws.Cells.Item(71, 11).Value = 30010
$ws.Cells.Item(71, 13).Value = -26266
$ws.Cells.Item(93, 8).Value = 2399.25
$ws.Cells.Item(93, 9).Value = 2399.25
$ws.Cells.Item(93, 11).Value = 2399.25
$ws.Cells.Item(93, 13).Value = -1151.25
$ws.Cells.Item(132, 8).Value = 1451665.4
$ws.Cells.Item(132, 9).Value = 1756854.4
$ws.Cells.Item(132, 10).Value = 2018
$ws.Cells.Item(132, 11).Value = 5270563.199999999
$ws.Cells.Item(132, 12).Value = 6054
$ws.Cells.Item(132, 13).Value = -5268033.199999999
$ws.Cells.Item(132, 14).Value = -11114

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 107650.414
$ws.Cells.Item(122, 9).Value = 12443.223
$ws.Cells.Item(122, 11).Value = 37329.669
$ws.Cells.Item(122, 13).Value = -34879.669
$ws.Cells.Item(132, 8).Value = 5210784.5
$ws.Cells.Item(132, 9).Value = 6175301.5
$ws.Cells.Item(132, 11).Value = 18525904.5
$ws.Cells.Item(132, 13).Value = -18523374.5
$ws.Cells.Item(136, 8).Value = 8833360
$ws.Cells.Item(136, 9).Value = 7248144
$ws.Cells.Item(136, 10).Value = 10418576
$ws.Cells.Item(136, 11).Value = 21744432
$ws.Cells.Item(136, 12).Value = 31255728
$ws.Cells.Item(136, 13).Value = -21741882
$ws.Cells.Item(136, 14).Value = -31260828